$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column C for rows 2-7
# from serial 45185 (2023-09-16) to serial 45204 (2023-10-05)
foreach ($row in 2..7) {
    $ws.Cells.Item($row, 3).Value = 45204
}
